# Apply the "plotted all the new graphs on one single scale" edit.
#
# Adds, on Sheet1, a second block of columns (P:AB) next to the existing
# leaf/spine data (C:N): a "no failure" header row mirroring the original
# leaf/spine headers, and a new "Probabilities" row (row 5) plus matching
# probability values for every data row (rows 7-21) that currently holds
# leaf/spine numbers. The probabilities are 0.5 in the first new column and
# 1/24 (as a shared formula) in the remaining eleven columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. New shared strings --------------------------------------------
# Order matters: the workbook's sharedStrings table gains "Probabilities"
# (index 14) before "no failure" (index 15), so we must write the B5 cell
# before the P3 cell.
$ws.Range("B5").Value = "Probabilities"
$ws.Range("P3").Value = "no failure"

# --- 2. Mirror the leaf/spine headers into the new block (row 3) ------
$headers = @("Leaf 1","Leaf 2","Leaf 3","Leaf 4","Leaf 5","Leaf 6","Leaf 7","Leaf 8","Spine 1","Spine 2","Spine 3","Spine 4")
$newCols = @("Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "3").Value = $headers[$i]
}

# --- 3. Probability values/formulas for row 5 and rows 7-21 -----------
$dataRows = @(5,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $dataRows) {
    $ws.Range("P$r").Value = 0.5
    $ws.Range("Q$r").Formula = "=1/24"
}

# One shared formula covering the whole rectangular block at once (mirrors
# how Excel groups a fill-across/fill-down into a single shared formula
# whose ref is the full filled rectangle, even though row 6 has no data).
$ws.Range("R5:AB21").Formula = "=1/24"

# Row 6 never had any data before this edit and must stay that way - remove
# whatever the block-fill above created there.
$ws.Range("P6:AB6").Clear()

# --- 4. Number formatting (percentage, 2 decimals) for the new values -
$ws.Range("P5:AB21").NumberFormat = "0.00%"

# Re-clear row 6 since applying the number format to the rectangular range
# re-creates (blank, styled-only) cells on that empty row.
$ws.Range("P6:AB6").Clear()

# --- 5. View state: scroll down a bit and leave the selection on Q21 --
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("Q21").Select()
